{"js": "// Update the date line and the 25 division-problem answers in the practice\n// table. Cells are addressed by absolute (row, column) position rather than\n// by searching for their old text, because a couple of the new answers are\n// identical to old answers that live elsewhere in the table (e.g. the new\n// value for row 0/col 3 \u2014 \"66\u00f76=11, 0\" \u2014 is the OLD value of row 16/col 1),\n// so a naive global text search-and-replace could touch the wrong cell.\n\n// 1) The date paragraph above the table.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(\"2025-04-15 Tuesday\", \"Replace\");\n\n// 2) The table of division problems. Only every 4th row (0, 4, 8, 12, 16)\n// actually holds text; the rows in between are blank spacer rows.\nconst table = body.tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst newValues = {\n  0: [\"30\u00f73=10, 0\", \"23\u00f76=3, 5\", \"77\u00f78=9, 5\", \"66\u00f76=11, 0\", \"13\u00f74=3, 1\"],\n  4: [\"66\u00f77=9, 3\", \"36\u00f74=9, 0\", \"35\u00f74=8, 3\", \"41\u00f74=10, 1\", \"40\u00f78=5, 0\"],\n  8: [\"54\u00f78=6, 6\", \"55\u00f75=11, 0\", \"20\u00f75=4, 0\", \"99\u00f79=11, 0\", \"41\u00f78=5, 1\"],\n  12: [\"30\u00f73=10, 0\", \"66\u00f72=33, 0\", \"54\u00f76=9, 0\", \"30\u00f78=3, 6\", \"14\u00f73=4, 2\"],\n  16: [\"30\u00f72=15, 0\", \"87\u00f79=9, 6\", \"10\u00f74=2, 2\", \"83\u00f76=13, 5\", \"28\u00f77=4, 0\"],\n};\n\nfor (const rowIndex of Object.keys(newValues)) {\n  const ri = Number(rowIndex);\n  const row = rows.items[ri];\n  row.cells.load(\"items\");\n  await context.sync();\n\n  const values = newValues[ri];\n  for (let c = 0; c < values.length; c++) {\n    row.cells.items[c].value = values[c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division-problem answers in the practice\n# table. Cells are addressed by absolute (row, column) position rather than\n# by searching for their old text, because a couple of the new answers are\n# identical to old answers that live elsewhere in the table (e.g. the new\n# value for row 1/col 4 -- \"66\u00f76=11, 0\" -- is the OLD value of row 17/col 2),\n# so a naive global Find/Replace could touch the wrong cell.\n\n$d = $word.ActiveDocument\n\n# 1) The date paragraph above the table.\n$d.Paragraphs.Item(1).Range.Text = \"2025-04-15 Tuesday\"\n\n# 2) The table of division problems. Only every 4th row (1, 5, 9, 13, 17 in\n# 1-based COM numbering) actually holds text; the rows in between are blank\n# spacer rows.\n$tbl = $d.Tables.Item(1)\n\n$tbl.Cell(1, 1).Range.Text = \"30\u00f73=10, 0\"\n$tbl.Cell(1, 2).Range.Text = \"23\u00f76=3, 5\"\n$tbl.Cell(1, 3).Range.Text = \"77\u00f78=9, 5\"\n$tbl.Cell(1, 4).Range.Text = \"66\u00f76=11, 0\"\n$tbl.Cell(1, 5).Range.Text = \"13\u00f74=3, 1\"\n\n$tbl.Cell(5, 1).Range.Text = \"66\u00f77=9, 3\"\n$tbl.Cell(5, 2).Range.Text = \"36\u00f74=9, 0\"\n$tbl.Cell(5, 3).Range.Text = \"35\u00f74=8, 3\"\n$tbl.Cell(5, 4).Range.Text = \"41\u00f74=10, 1\"\n$tbl.Cell(5, 5).Range.Text = \"40\u00f78=5, 0\"\n\n$tbl.Cell(9, 1).Range.Text = \"54\u00f78=6, 6\"\n$tbl.Cell(9, 2).Range.Text = \"55\u00f75=11, 0\"\n$tbl.Cell(9, 3).Range.Text = \"20\u00f75=4, 0\"\n$tbl.Cell(9, 4).Range.Text = \"99\u00f79=11, 0\"\n$tbl.Cell(9, 5).Range.Text = \"41\u00f78=5, 1\"\n\n$tbl.Cell(13, 1).Range.Text = \"30\u00f73=10, 0\"\n$tbl.Cell(13, 2).Range.Text = \"66\u00f72=33, 0\"\n$tbl.Cell(13, 3).Range.Text = \"54\u00f76=9, 0\"\n$tbl.Cell(13, 4).Range.Text = \"30\u00f78=3, 6\"\n$tbl.Cell(13, 5).Range.Text = \"14\u00f73=4, 2\"\n\n$tbl.Cell(17, 1).Range.Text = \"30\u00f72=15, 0\"\n$tbl.Cell(17, 2).Range.Text = \"87\u00f79=9, 6\"\n$tbl.Cell(17, 3).Range.Text = \"10\u00f74=2, 2\"\n$tbl.Cell(17, 4).Range.Text = \"83\u00f76=13, 5\"\n$tbl.Cell(17, 5).Range.Text = \"28\u00f77=4, 0\"\n"}
